$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BE")

# Copy the formatting of the prior data row (17) down into the new row (19)
# so the new row's styles (incl. the date number format on column G) match
# the rest of the table.
$ws.Range("A17:I17").Copy()
$ws.Range("A19:I19").PasteSpecial(-4122)  # xlPasteFormats

# New row 19: BE Wave10 entry, continuing the existing "BE" sheet table.
$ws.Range("A19").Value = "be"
$ws.Range("B19").Value = 5
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = "B"
$ws.Range("F19").Value = 10
$ws.Range("G19").Value = 44285
$ws.Range("H19").Value = "21-019042_BE_Wave10_Final_v1_260321_IntClientUse"

# Extend column I's formula down into the new row - same pattern as every
# other row in the table (relative references shift per-row automatically).
$ws.Range("I19").Formula = '=A19&"_"&"sr"&TEXT(D19,"00")&"_"&YEAR(G19)&TEXT(G19,"MM")&TEXT(G19,"DD")&"_p"&E19&"_wv"&TEXT(F19,"00")&""'

$ws.Range("B19").Select()
